# fix: modified the cell name of excel template
# The "Kategorie" header in I1 is renamed to "Kunden Interessen".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the I1 header cell (was "Kategorie") to its new label.
$ws.Range("I1").Value = "Kunden Interessen"

# New column I is now wide enough to show the longer header text without
# being auto-fit from the old short label.
$ws.Columns.Item(9).ColumnWidth = 21.86

# The active selection in the saved file moved to J5.
$null = $ws.Range("J5").Select()
